$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 59.87124934497245
$ws.Range("C2").Value = 0.2959971687491965
$ws.Range("D2").Value = 1669
$ws.Range("E2").Value = 59.29068453228387
$ws.Range("F2").Value = 60.45181415766104

$ws.Range("B3").Value = 62.07504175816011
$ws.Range("C3").Value = 0.3073559955958241
$ws.Range("D3").Value = 1669
$ws.Range("E3").Value = 61.47219789730211
$ws.Range("F3").Value = 62.67788561901811

$ws.Range("B4").Value = 63.09915708230246
$ws.Range("C4").Value = 0.6420788306891397
$ws.Range("D4").Value = 1669
$ws.Range("E4").Value = 61.83979241457418
$ws.Range("F4").Value = 64.35852175003075

$ws.Range("B5").Value = 111.4734215176497
$ws.Range("C5").Value = 0.5216112834330607
$ws.Range("D5").Value = 1669
$ws.Range("E5").Value = 110.4503402552774
$ws.Range("F5").Value = 112.4965027800221

$ws.Range("B6").Value = 110.8226110978106
$ws.Range("C6").Value = 0.5416280027645338
$ws.Range("D6").Value = 1669
$ws.Range("E6").Value = 109.760269314946
$ws.Range("F6").Value = 111.8849528806752

$ws.Range("B7").Value = 108.7216054737588
$ws.Range("C7").Value = 1.131482319091845
$ws.Range("D7").Value = 1669
$ws.Range("E7").Value = 106.5023314738761
$ws.Range("F7").Value = 110.9408794736416

$ws.Range("B8").Value = 93.63818078330706
$ws.Range("C8").Value = 0.7238575381130589
$ws.Range("D8").Value = 1669
$ws.Range("E8").Value = 92.21841647338587
$ws.Range("F8").Value = 95.05794509322826

$ws.Range("B9").Value = 101.3472329196135
$ws.Range("C9").Value = 0.751635413394086
$ws.Range("D9").Value = 1669
$ws.Range("E9").Value = 99.87298546369611
$ws.Range("F9").Value = 102.8214803755308

$ws.Range("B10").Value = 103.828785911391
$ws.Range("C10").Value = 1.570196105662629
$ws.Range("D10").Value = 1669
$ws.Range("E10").Value = 100.7490246696854
$ws.Range("F10").Value = 106.9085471530966

$ws.Range("B11").Value = 18.02789681008836
$ws.Range("C11").Value = 0.1335004467845181
$ws.Range("D11").Value = 1669
$ws.Range("E11").Value = 17.76605085322055
$ws.Range("F11").Value = 18.28974276695617

$ws.Range("B12").Value = 18.97275615880135
$ws.Range("C12").Value = 0.1386234973372673
$ws.Range("D12").Value = 1669
$ws.Range("E12").Value = 18.70086192040058
$ws.Range("F12").Value = 19.24465039720212

$ws.Range("B13").Value = 19.03805082755149
$ws.Range("C13").Value = 0.289589968478757
$ws.Range("D13").Value = 1669
$ws.Range("E13").Value = 18.47005301015068
$ws.Range("F13").Value = 19.6060486449523

$ws.Range("B14").Value = 44.19184981096711
$ws.Range("C14").Value = 0.2158068519162584
$ws.Range("D14").Value = 1669
$ws.Range("E14").Value = 43.76856919281418
$ws.Range("F14").Value = 44.61513042912003

$ws.Range("B15").Value = 44.78081080015667
$ws.Range("C15").Value = 0.2240883928296095
$ws.Range("D15").Value = 1669
$ws.Range("E15").Value = 44.34128688052404
$ws.Range("F15").Value = 45.22033471978929

$ws.Range("B16").Value = 44.91683844128804
$ws.Range("C16").Value = 0.4681295152877089
$ws.Range("D16").Value = 1669
$ws.Range("E16").Value = 43.99865559006189
$ws.Range("F16").Value = 45.8350212925142

$ws.Range("B17").Value = 115.2833173892807
$ws.Range("C17").Value = 0.469993197060413
$ws.Range("D17").Value = 1669
$ws.Range("E17").Value = 114.3614791380254
$ws.Range("F17").Value = 116.2051556405361

$ws.Range("B18").Value = 116.5341795233623
$ws.Range("C18").Value = 0.4880290835760223
$ws.Range("D18").Value = 1669
$ws.Range("E18").Value = 115.5769659301033
$ws.Range("F18").Value = 117.4913931166212

$ws.Range("B19").Value = 114.5415281333033
$ws.Range("C19").Value = 1.019512057076795
$ws.Range("D19").Value = 1669
$ws.Range("E19").Value = 112.5418710793128
$ws.Range("F19").Value = 116.5411851872938

$ws.Range("B20").Value = 148.9833993131405
$ws.Range("C20").Value = 0.6253200781321092
$ws.Range("D20").Value = 1669
$ws.Range("E20").Value = 147.7569050342753
$ws.Range("F20").Value = 150.2098935920056

$ws.Range("B21").Value = 155.4663464356074
$ws.Range("C21").Value = 0.6493165998597905
$ws.Range("D21").Value = 1669
$ws.Range("E21").Value = 154.1927857060702
$ws.Range("F21").Value = 156.7399071651447

$ws.Range("B22").Value = 156.9495458206324
$ws.Range("C22").Value = 1.35644805749378
$ws.Range("D22").Value = 1669
$ws.Range("E22").Value = 154.2890270874686
$ws.Range("F22").Value = 159.6100645537961

$ws.Range("B23").Value = 280.549228648549
$ws.Range("C23").Value = 0.663242606171626
$ws.Range("D23").Value = 1669
$ws.Range("E23").Value = 279.248353640027
$ws.Range("F23").Value = 281.8501036570711

$ws.Range("B24").Value = 282.1124659539652
$ws.Range("C24").Value = 0.6886943966486927
$ws.Range("D24").Value = 1669
$ws.Range("E24").Value = 280.7616701504796
$ws.Range("F24").Value = 283.4632617574508

$ws.Range("B25").Value = 286.1888316112021
$ws.Range("C25").Value = 1.438709832372514
$ws.Range("D25").Value = 1669
$ws.Range("E25").Value = 283.3669657539232
$ws.Range("F25").Value = 289.010697468481

$ws.Range("B26").Value = 146.3033083386354
$ws.Range("C26").Value = 0.3285988775409958
$ws.Range("D26").Value = 1669
$ws.Range("E26").Value = 145.6587989786458
$ws.Range("F26").Value = 146.9478176986249

$ws.Range("B27").Value = 147.1214228131406
$ws.Range("C27").Value = 0.3412087878578981
$ws.Range("D27").Value = 1669
$ws.Range("E27").Value = 146.4521805469138
$ws.Range("F27").Value = 147.7906650793673

$ws.Range("B28").Value = 148.3230816569054
$ws.Range("C28").Value = 0.7127986525981227
$ws.Range("D28").Value = 1669
$ws.Range("E28").Value = 146.9250080943135
$ws.Range("F28").Value = 149.7211552194973

